$p = $ppt.ActivePresentation
$p.Designs.Item(1).SlideMaster.ApplyTheme("C:\temp\integral.xml")
